$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "VALOR MORA" total value (E11): 66919 -> 225210 ---
$ws.Range("E11").Value = 225210

# --- 2. Update worker / period counts (C13, F13): 2 -> 3 ---
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 3

# --- 3. Insert two new data rows right after the existing data row 17 ---
#     (old rows 18-21 were blank, old rows 22/23 -> shift to 24/25)
$ws.Range("18:19").EntireRow.Insert()

# --- 4. Preserve the "last row" (total) border styling: it currently still
#     sits on row 17 (pre-insert last data row) - copy that formatting
#     down onto the new final data row (19) before row 17 gets restyled. ---
$ws.Range("B17:J17").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# --- 5. Row 17 is no longer the last row, so it now takes the regular
#     (non-total) row styling, same as row 16. Also stamp that same
#     regular styling onto the brand-new row 18. ---
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$ws.Range("B18:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 6. New row 18 data: MARIELE DEL MAR PIÑA PAJARO, doc 1143401993, period 2508 ---
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143401993"
$ws.Range("D18").Value = "MARIELE DEL MAR PIÑA PAJARO"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 36051
$ws.Range("G18").Value = 2458000

# --- 7. New row 19 data: LINDA YERALDIN NIETO ESTRADA, doc 1065827176, new
#     period 2508 (she already had period 2507 on row 17) ---
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1065827176"
$ws.Range("D19").Value = "LINDA YERALDIN NIETO ESTRADA"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 122240
$ws.Range("G19").Value = 3056000
